$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("Infusion-related reactions" / "Bullet point"), shifting rows 3-8 up.
$ws.Rows("2:2").Delete()
